# Updated symbol list on Wed Dec 28 15:17:27 UTC 2022 with GitHub Actions
#
# Refreshes the crypto price/volume snapshot on Sheet1: the "Price" column
# (D) gets new quotes and a handful of "Volume(1h)" rank labels (E) flip
# between "...Bestin24h" / "...Worstin24h" / plain as the 24h leaders churn.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates: cell -> new text. These are stored as plain text
# in the sheet (not numbers), so force Text format before writing the value
# to keep "23.90"-style strings from being auto-converted to numerics.
$priceUpdates = [ordered]@{
    "D2"  = "244.55"
    "D3"  = "23.87"
    "D4"  = "5.294"
    "D5"  = "0.05890"
    "D6"  = "6.478"
    "D7"  = "3.327"
    "D8"  = "0.8168"
    "D9"  = "0.8926"
    "D11" = "0.07239"
    "D12" = "0.03083"
    "D13" = "0.03033"
    "D14" = "0.09360"
    "D15" = "3.855"
    "D16" = "0.001541"
    "D17" = "0.04701"
    "D18" = "0.0005998"
    "D19" = "0.006167"
    "D21" = "0.004615"
    "D23" = "3.561"
    "D24" = "2.175"
    "D25" = "0.3199"
    "D26" = "0.1314"
    "D28" = "0.0002339"
    "D40" = "0.03809"
    "D41" = "0.006354"
    "D42" = "0.1063"
    "D43" = "0.002501"
    "D44" = "0.007089"
    "D45" = "0.00005393"
    "D46" = "0.00000000750"
    "D47" = "0.5398"
    "D48" = "0.002028"
    "D49" = "0.00002099"
    "D50" = "0.0001999"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
}

# Volume(1h) column (E) rank-label updates (plain text, no formatting needed).
$volumeUpdates = [ordered]@{
    "E19" = "18TigerCashTCHBestin24h"
    "E43" = "42CEJICEJIWorstin24h"
    "E47" = "46CoinbaseStockTokenCOIN"
    "E48" = "47BOLOBOLO"
}

foreach ($cellRef in $volumeUpdates.Keys) {
    $ws.Range($cellRef).Value = $volumeUpdates[$cellRef]
}
